$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J2").Value = 3.8
$ws.Range("L2").Value = 1.34
$ws.Range("Q2").Value = 1.72
$ws.Range("S2").Value = 2.82
$ws.Range("X2").Value = 21
$ws.Range("AO2").Value = 27
$ws.Range("F3").Value = 1.53
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = 8.199999999999999
$ws.Range("M3").Value = 1.05
$ws.Range("N3").Value = 3.6
$ws.Range("P3").Value = 1.84
$ws.Range("R3").Value = 1.34
$ws.Range("S3").Value = 3.35
$ws.Range("T3").Value = 1.98
$ws.Range("U3").Value = 1.84
$ws.Range("N4").Value = 2.96
$ws.Range("Q4").Value = 2.28
$ws.Range("T4").Value = 1.98
$ws.Range("W4").Value = 1.92
$ws.Range("X4").Value = 12.5
$ws.Range("AD4").Value = 20
$ws.Range("AF4").Value = 12.5
$ws.Range("AK4").Value = 980
$ws.Range("AL4").Value = 50
$ws.Range("G5").Value = 4.6
$ws.Range("H5").Value = 1.93
$ws.Range("I5").Value = 2.12
$ws.Range("K5").Value = 4.2
$ws.Range("N5").Value = 3.85
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 1.98
$ws.Range("Q5").Value = 1.84
$ws.Range("S5").Value = 3.1
$ws.Range("T5").Value = 1.74
$ws.Range("U5").Value = 2.12
$ws.Range("V5").Value = 1.89
$ws.Range("W5").Value = 1.28
$ws.Range("AF5").Value = 38
$ws.Range("AI5").Value = 42
$ws.Range("F6").Value = 8.6
$ws.Range("G6").Value = 9.6
$ws.Range("N6").Value = 4.8
$ws.Range("P6").Value = 2.32
$ws.Range("Q6").Value = 1.68
$ws.Range("R6").Value = 1.53
$ws.Range("U6").Value = 1.93
$ws.Range("Y6").Value = 11.5
$ws.Range("AL6").Value = 130
$ws.Range("AM6").Value = 160
$ws.Range("AN6").Value = 190
$ws.Range("I7").Value = 21
$ws.Range("R7").Value = 1.95
$ws.Range("V7").Value = 1.05
$ws.Range("AB7").Value = 17
$ws.Range("AF7").Value = 12.5
$ws.Range("AN7").Value = 2.94
$ws.Range("F8").Value = 1.3
$ws.Range("I8").Value = 10.5
$ws.Range("K8").Value = 7.6
$ws.Range("T8").Value = 1.79
$ws.Range("U8").Value = 2.04
$ws.Range("AA8").Value = 440
$ws.Range("AB8").Value = 19.5
$ws.Range("AE8").Value = 160
$ws.Range("AO8").Value = 120
$ws.Range("H9").Value = 2.84
$ws.Range("S9").Value = 3.75
$ws.Range("T9").Value = 1.6
$ws.Range("G10").Value = 3.85
$ws.Range("I10").Value = 2.62
$ws.Range("J10").Value = 2.94
$ws.Range("K10").Value = 3.6
$ws.Range("M10").Value = 1.08
$ws.Range("N10").Value = 3.4
$ws.Range("O10").Value = 1.36
$ws.Range("P10").Value = 1.82
$ws.Range("Q10").Value = 1.92
$ws.Range("R10").Value = 1.32
$ws.Range("S10").Value = 3.7
$ws.Range("T10").Value = 1.8
$ws.Range("U10").Value = 2.06
$ws.Range("V10").Value = 1.6
$ws.Range("W10").Value = 1.35
$ws.Range("X10").Value = 1000
$ws.Range("Y10").Value = 1000
$ws.Range("Z10").Value = 1000
$ws.Range("AC10").Value = 1000
$ws.Range("AG10").Value = 1000
$ws.Range("AH10").Value = 1000
$ws.Range("AI10").Value = 980
$ws.Range("AJ10").Value = 80
$ws.Range("AK10").Value = 980
$ws.Range("AL10").Value = 980
$ws.Range("AM10").Value = 120
$ws.Range("AN10").Value = 980
$ws.Range("AO10").Value = 980
$ws.Range("I11").Value = 4.8
$ws.Range("J11").Value = 3.7
$ws.Range("R11").Value = 1.14
$ws.Range("S11").Value = 1.75
$ws.Range("T11").Value = 1.47
$ws.Range("U11").Value = 1.83
$ws.Range("V11").Value = 1.26
$ws.Range("Q12").Value = 1.32
$ws.Range("R12").Value = 1.86
$ws.Range("S12").Value = 1.92
$ws.Range("AA12").Value = 80
$ws.Range("AG12").Value = 14
$ws.Range("F13").Value = 2.36
$ws.Range("I13").Value = 2.94
$ws.Range("K13").Value = 5
$ws.Range("L14").Value = 1.24
$ws.Range("G15").Value = 1.87
$ws.Range("H15").Value = 2.2
$ws.Range("F19").Value = 1.67
$ws.Range("G19").Value = 1.95
$ws.Range("K19").Value = 5.3
$ws.Range("W19").Value = 2.06
$ws.Range("L20").Value = 1.57
$ws.Range("N20").Value = 1.03
$ws.Range("O20").Value = 1.54
$ws.Range("Q20").Value = 1.01
$ws.Range("P21").Value = 1.81
$ws.Range("Q21").Value = 2.1
$ws.Range("R21").Value = 1.29
$ws.Range("V21").Value = 1.24
$ws.Range("F22").Value = 2.62
$ws.Range("N22").Value = 2.8
$ws.Range("AE22").Value = 44
$ws.Range("AH22").Value = 22
$ws.Range("AN22").Value = 40
$ws.Range("AO22").Value = 50
$ws.Range("O23").Value = 1.54
$ws.Range("P24").Value = 1.35
$ws.Range("U24").Value = 1.57
$ws.Range("R25").Value = 1.27
$ws.Range("AH26").Value = 29
$ws.Range("H27").Value = 10
$ws.Range("V27").Value = 1.1
$ws.Range("AD27").Value = 42
$ws.Range("AE27").Value = 260
$ws.Range("AM27").Value = 300
$ws.Range("N28").Value = 1.03
$ws.Range("O28").Value = 1.47
$ws.Range("P28").Value = 1.08
$ws.Range("Q28").Value = 1.47
$ws.Range("S28").Value = 1.01
